$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

# New row 53: 2020-07-22 raw/clean data.
# Force column A to be stored as text (shared string) rather than
# letting Excel auto-detect the "YYYY-MM-DD" literal as a date, then
# drop the temporary number format back to the sheet's default style
# so no extra formatting is left behind on the cell.
$row = 53
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2020-07-22"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = 362274
$ws.Range("C$row").Value = 411673
$ws.Range("D$row").Value = 87905
$ws.Range("E$row").Value = 41190
$ws.Range("F$row").Value = 28.24
